$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list with latest scraped price/volume data (GitHub Actions run)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.899.60"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.74%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.798.97"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.45%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.006"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.38%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.16"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.67%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4980"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.96%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3848"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.16%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09340"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +16.76%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.091"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.17%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.52"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.12%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.298"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.42%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.006"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.35%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.64"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.87%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.807.86"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.10%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.180"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.41%  "

# Row 17
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.98"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.70%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001106"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.31%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06589"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.54%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.34%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.07"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.931"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.22%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.977.46"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.61%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.97"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.233"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.73%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.57"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.012.76"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.98%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.37"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.16%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.371"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.90%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.52"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.08%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1072"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.55%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.046"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.58%  "

# Row 33
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.624"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.29%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.529"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06817"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.832"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.08%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02294"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.41%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2131"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.58%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.31"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -6.31%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.928"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.47%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6117"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.57%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.50%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.144"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.55%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.04"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.82%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5846"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.58%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.281"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.86%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.661"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.69%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.65"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.99%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.939"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.21%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.169"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.68%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06722"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.50%  "
